# Uruguay Primera Division - odds base update (19-04-2024 21:40)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Rows 118 and 120 (match ids 7013886 / 7013409, identical kickoff timestamp)
# swap their complete data (everything except the leading row-index column A).
# ---------------------------------------------------------------------------

# New content for row 118 (was row 120's content): Racing Club de Montevideo vs Cerro
$ws.Cells.Item(118, 2).Value = 7013886          # B118
$ws.Cells.Item(118, 6).Value = "Racing Club de Montevideo"   # F118
$ws.Cells.Item(118, 7).Value = "Cerro"                        # G118
$ws.Cells.Item(118, 8).Value = 0                # H118
$ws.Cells.Item(118, 9).Value = 1                # I118
$ws.Cells.Item(118, 10).Value = "A"             # J118
$ws.Cells.Item(118, 11).Value = 2.25            # K118
$ws.Cells.Item(118, 12).Value = 3.1             # L118
$ws.Cells.Item(118, 13).Value = 3.25            # M118
$ws.Cells.Item(118, 14).Value = 2.25            # N118
$ws.Cells.Item(118, 15).Value = 2.875           # O118
$ws.Cells.Item(118, 16).Value = 3.5             # P118
$ws.Cells.Item(118, 17).Value = -0.25           # Q118
$ws.Cells.Item(118, 18).Value = 1.95            # R118
$ws.Cells.Item(118, 19).Value = 1.9             # S118
$ws.Cells.Item(118, 20).Value = 2               # T118
$ws.Cells.Item(118, 21).Value = 1.925           # U118
$ws.Cells.Item(118, 22).Value = 1.925           # V118
$ws.Cells.Item(118, 23).Value = -1              # W118
$ws.Cells.Item(118, 24).Value = -1              # X118
$ws.Cells.Item(118, 25).Value = 2.5             # Y118
$ws.Cells.Item(118, 26).Value = -1              # Z118
$ws.Cells.Item(118, 27).Value = 0.8999999999999999   # AA118
$ws.Cells.Item(118, 28).Value = -1              # AB118
$ws.Cells.Item(118, 29).Value = 0.925           # AC118

# New content for row 120 (was row 118's content): Nacional De Football vs Torque
$ws.Cells.Item(120, 2).Value = 7013409          # B120
$ws.Cells.Item(120, 6).Value = "Nacional De Football"        # F120
$ws.Cells.Item(120, 7).Value = "Torque"                       # G120
$ws.Cells.Item(120, 8).Value = 1                # H120
$ws.Cells.Item(120, 9).Value = 1                # I120
$ws.Cells.Item(120, 10).Value = "D"             # J120
$ws.Cells.Item(120, 11).Value = 1.666           # K120
$ws.Cells.Item(120, 12).Value = 3.9             # L120
$ws.Cells.Item(120, 13).Value = 4.5             # M120
$ws.Cells.Item(120, 14).Value = 1.615           # N120
$ws.Cells.Item(120, 15).Value = 4               # O120
$ws.Cells.Item(120, 16).Value = 4.75            # P120
$ws.Cells.Item(120, 17).Value = -0.75           # Q120
$ws.Cells.Item(120, 18).Value = 1.8             # R120
$ws.Cells.Item(120, 19).Value = 2.05            # S120
$ws.Cells.Item(120, 20).Value = 2.75            # T120
$ws.Cells.Item(120, 21).Value = 1.95            # U120
$ws.Cells.Item(120, 22).Value = 1.9             # V120
$ws.Cells.Item(120, 23).Value = -1              # W120
$ws.Cells.Item(120, 24).Value = 3               # X120
$ws.Cells.Item(120, 25).Value = -1              # Y120
$ws.Cells.Item(120, 26).Value = -1              # Z120
$ws.Cells.Item(120, 27).Value = 1.05            # AA120
$ws.Cells.Item(120, 28).Value = -1              # AB120
$ws.Cells.Item(120, 29).Value = 0.8999999999999999   # AC120

# ---------------------------------------------------------------------------
# Rows 185-187: refreshed odds for upcoming fixtures (values shifted up from
# the old rows 187-189, with updated odds), then drop the now-obsolete
# trailing rows 188-190 (one fixture removed entirely).
# ---------------------------------------------------------------------------

# Row 185 (was row 187's fixture: Danubio vs Cerro Largo)
$ws.Cells.Item(185, 2).Value = 8081162
$ws.Cells.Item(185, 5).Value = 45402.41666666666
$ws.Cells.Item(185, 6).Value = "Danubio"
$ws.Cells.Item(185, 7).Value = "Cerro Largo"
$ws.Cells.Item(185, 11).Value = 2.3
$ws.Cells.Item(185, 12).Value = 3
$ws.Cells.Item(185, 13).Value = 3.4
$ws.Cells.Item(185, 14).Value = 2.3
$ws.Cells.Item(185, 15).Value = 3
$ws.Cells.Item(185, 16).Value = 3.4
$ws.Cells.Item(185, 17).Value = -0.25
$ws.Cells.Item(185, 18).Value = 1.975
$ws.Cells.Item(185, 19).Value = 1.875
$ws.Cells.Item(185, 20).Value = 2
$ws.Cells.Item(185, 21).Value = 1.975
$ws.Cells.Item(185, 22).Value = 1.875
$ws.Cells.Item(185, 23).Value = 0
$ws.Cells.Item(185, 24).Value = 0
$ws.Cells.Item(185, 25).Value = 0
$ws.Cells.Item(185, 26).Value = 0
$ws.Cells.Item(185, 27).Value = 0

# Row 186 (was row 188's fixture: Boston River vs Penarol)
$ws.Cells.Item(186, 2).Value = 8081144
$ws.Cells.Item(186, 5).Value = 45402.625
$ws.Cells.Item(186, 6).Value = "Boston River"
$ws.Cells.Item(186, 7).Value = "Penarol"
$ws.Cells.Item(186, 11).Value = 4.75
$ws.Cells.Item(186, 12).Value = 3.75
$ws.Cells.Item(186, 13).Value = 1.727
$ws.Cells.Item(186, 14).Value = 4
$ws.Cells.Item(186, 15).Value = 3.6
$ws.Cells.Item(186, 16).Value = 1.909
$ws.Cells.Item(186, 17).Value = 0.5
$ws.Cells.Item(186, 18).Value = 1.925
$ws.Cells.Item(186, 19).Value = 1.925
$ws.Cells.Item(186, 20).Value = 2.25
$ws.Cells.Item(186, 21).Value = 1.95
$ws.Cells.Item(186, 22).Value = 1.9
$ws.Cells.Item(186, 23).Value = 0
$ws.Cells.Item(186, 24).Value = 0
$ws.Cells.Item(186, 25).Value = 0
$ws.Cells.Item(186, 26).Value = 0
$ws.Cells.Item(186, 27).Value = 0

# Row 187 (was row 189's fixture: Nacional De Football vs Rampla Juniors)
$ws.Cells.Item(187, 2).Value = 8081249
$ws.Cells.Item(187, 5).Value = 45402.75
$ws.Cells.Item(187, 6).Value = "Nacional De Football"
$ws.Cells.Item(187, 7).Value = "Rampla Juniors"
$ws.Cells.Item(187, 11).Value = 1.444
$ws.Cells.Item(187, 12).Value = 4
$ws.Cells.Item(187, 13).Value = 8.5
$ws.Cells.Item(187, 14).Value = 1.25
$ws.Cells.Item(187, 15).Value = 5
$ws.Cells.Item(187, 16).Value = 13
$ws.Cells.Item(187, 17).Value = -1.5
$ws.Cells.Item(187, 18).Value = 1.85
$ws.Cells.Item(187, 19).Value = 2
$ws.Cells.Item(187, 20).Value = 2.5
$ws.Cells.Item(187, 21).Value = 1.875
$ws.Cells.Item(187, 22).Value = 1.975
$ws.Cells.Item(187, 23).Value = 0
$ws.Cells.Item(187, 24).Value = 0
$ws.Cells.Item(187, 25).Value = 0
$ws.Cells.Item(187, 26).Value = 0
$ws.Cells.Item(187, 27).Value = 0

# Remove the now-obsolete trailing rows 188:190 (their fixtures were dropped
# from the feed / folded into the refreshed rows above).
$ws.Range("A188:AC190").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
